# Updates cryptos list: Coin/Link swaps (rows 27-28, 35, 37, 48-49) and refreshed
# Price (D) / Volume(1h) (E) figures, matching the Nov 9 2023 GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.726.85"
$ws.Range("E2").Value = "  +2.57%  "

$ws.Range("D3").Value = "'2.127.62"
$ws.Range("E3").Value = "  +12.47%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'250.98"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("E6").Value = "  -2.78%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'45.37"
$ws.Range("E8").Value = "  +4.66%  "

$ws.Range("D9").Value = "'61.27"
$ws.Range("E9").Value = "  +7.22%  "

$ws.Range("D10").Value = "'0.369"
$ws.Range("E10").Value = "  +2.78%  "

$ws.Range("D11").Value = "'0.0736"
$ws.Range("E11").Value = "  -2.88%  "

$ws.Range("E12").Value = "  +0.59%  "

$ws.Range("D13").Value = "'14.69"
$ws.Range("E13").Value = "  -1.70%  "

$ws.Range("D14").Value = "'2.419.76"
$ws.Range("E14").Value = "  +11.54%  "

$ws.Range("D15").Value = "'0.852"
$ws.Range("E15").Value = "  +7.84%  "

$ws.Range("D16").Value = "'2.112.94"
$ws.Range("E16").Value = "  +11.39%  "

$ws.Range("D17").Value = "'5.10"
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("D18").Value = "'36.836.42"
$ws.Range("E18").Value = "  +2.45%  "

$ws.Range("D19").Value = "'73.56"
$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("D20").Value = "'0.0₃0823"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").Value = "'241.44"
$ws.Range("E21").Value = "  -2.18%  "

$ws.Range("D22").Value = "'12.99"
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").Value = "'5.14"
$ws.Range("E23").Value = "  -0.95%  "

$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  -8.76%  "

$ws.Range("D26").Value = "'170.04"
$ws.Range("E26").Value = "  +1.71%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'20.74"
$ws.Range("E27").Value = "  +12.41%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'9.05"
$ws.Range("E28").Value = "  +4.03%  "

$ws.Range("E29").Value = "  -7.74%  "

$ws.Range("E30").Value = "  -3.89%  "

$ws.Range("D31").Value = "'22.08"
$ws.Range("E31").Value = "  +44.93%  "

$ws.Range("D32").Value = "'4.48"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").Value = "'0.0599"
$ws.Range("E33").Value = "  -1.59%  "

$ws.Range("D34").Value = "'0.0910"
$ws.Range("E34").Value = "  +16.97%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'2.34"
$ws.Range("E35").Value = "  +19.85%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.88"
$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("D38").Value = "'4.09"
$ws.Range("E38").Value = "  -4.49%  "

$ws.Range("D39").Value = "'0.908"
$ws.Range("E39").Value = "  +5.43%  "

$ws.Range("E40").Value = "  -8.72%  "

$ws.Range("D41").Value = "'1.21"
$ws.Range("E41").Value = "  +10.48%  "

$ws.Range("D42").Value = "'100.59"
$ws.Range("E42").Value = "  +0.67%  "

$ws.Range("D43").Value = "'0.0220"
$ws.Range("E43").Value = "  -3.46%  "

$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  +16.49%  "

$ws.Range("D45").Value = "'16.35"
$ws.Range("E45").Value = "  -3.80%  "

$ws.Range("D46").Value = "'1.371.18"
$ws.Range("E46").Value = "  +4.09%  "

$ws.Range("D47").Value = "'0.0842"
$ws.Range("E47").Value = "  +4.12%  "

$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "'2.300.60"
$ws.Range("E48").Value = "  +11.10%  "

$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.84"
$ws.Range("E49").Value = "  +3.16%  "

$ws.Range("D50").Value = "'2.28"
$ws.Range("E50").Value = "  -2.74%  "

$ws.Range("D51").Value = "'3.93"
$ws.Range("E51").Value = "  +17.31%  "
